$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (date number format / style) from the row above into the new row's date cell
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)

# Populate the new row of daily expense data
$ws.Range("A14").Value = 43802
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 2.5
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 3

# Match the updated selection/active cell recorded in the workbook
$ws.Range("K14").Select()
